$wb = $excel.ActiveWorkbook
$wsLinear = $wb.Worksheets.Item("Simple Linear Regression")
$wsRPart = $wb.Worksheets.Item("RPart")

# --- RPart sheet: new header labels (shared strings reshuffle) ---
$wsRPart.Range("B1").Value = "TrainCP0.01"
$wsRPart.Range("C1").Value = "TestCP0.01"
$wsRPart.Range("D1").Value = "TrainCP0.001"
$wsRPart.Range("E1").Value = "TestCP0.001"
$wsRPart.Range("F1").Value = "TrainCP0.0001"
$wsRPart.Range("G1").Value = "TestCP0.0001"
$wsRPart.Range("H1").Value = "TrainCP0.0001_MinSp10"
$wsRPart.Range("I1").Value = "TrainCP0.0001_MinSp10.1"
$wsRPart.Range("J1").Value = "Train0.00001_MinSp10"
$wsRPart.Range("K1").Value = "Test0.00001_MinSp10"
$wsRPart.Range("L1").Value = "Train0.00002_MinSp10"
$wsRPart.Range("M1").Value = "Test0.00002_MinSp10"

# --- RPart sheet: refreshed numeric data (rows 2-5) ---
$wsRPart.Range("B2").Value = 1.05569027399669
$wsRPart.Range("C2").Value = 1.05190514749282
$wsRPart.Range("D2").Value = 0.81217515351616598
$wsRPart.Range("E2").Value = 0.80809969373582002
$wsRPart.Range("F2").Value = 0.68989815808176602
$wsRPart.Range("G2").Value = 0.70770995155222904
$wsRPart.Range("H2").Value = 0.68552576753924599
$wsRPart.Range("I2").Value = 0.71067366279937405
$wsRPart.Range("J2").Value = 0.454464456272199
$wsRPart.Range("K2").Value = 0.80405455103368895
$wsRPart.Range("L2").Value = 0.72946724337832602
$wsRPart.Range("M2").Value = 0.73368388849503097

$wsRPart.Range("B3").Value = 1.91940610862456
$wsRPart.Range("C3").Value = 1.8989350078844101
$wsRPart.Range("D3").Value = 1.2200367332072899
$wsRPart.Range("E3").Value = 1.2303173328936401
$wsRPart.Range("F3").Value = 0.92771673787272002
$wsRPart.Range("G3").Value = 1.02770054103641
$wsRPart.Range("H3").Value = 0.90739142363824299
$wsRPart.Range("I3").Value = 1.0400438064044799
$wsRPart.Range("J3").Value = 0.42847650866287101
$wsRPart.Range("K3").Value = 1.3991635617348801
$wsRPart.Range("L3").Value = 1.01868441349299
$wsRPart.Range("M3").Value = 1.0692727891987499

$wsRPart.Range("B4").Value = 1.38542632738972
$wsRPart.Range("C4").Value = 1.37801850781635
$wsRPart.Range("D4").Value = 1.10455272993519
$wsRPart.Range("E4").Value = 1.1091967061318
$wsRPart.Range("F4").Value = 0.96318053233686196
$wsRPart.Range("G4").Value = 1.01375566140782
$wsRPart.Range("H4").Value = 0.952570954647602
$wsRPart.Range("I4").Value = 1.0198253803492501
$wsRPart.Range("J4").Value = 0.65458117041576402
$wsRPart.Range("K4").Value = 1.1828624441307101
$wsRPart.Range("L4").Value = 1.00929897131276
$wsRPart.Range("M4").Value = 1.0340564729253201

$wsRPart.Range("B5").Value = 0.31414610963362999
$wsRPart.Range("C5").Value = 0.31886647906499699
$wsRPart.Range("D5").Value = 0.247710487123041
$wsRPart.Range("E5").Value = 0.24863476468752399
$wsRPart.Range("F5").Value = 0.188502275545435
$wsRPart.Range("G5").Value = 0.18948469970419701
$wsRPart.Range("H5").Value = 0.188092578765611
$wsRPart.Range("I5").Value = 0.189911559401402
$wsRPart.Range("J5").Value = 0.118854461751733
$wsRPart.Range("K5").Value = 0.209367742733573
$wsRPart.Range("L5").Value = 0.19473742295281701
$wsRPart.Range("M5").Value = 0.19392897955701499

# --- Highlight the TrainCP0.0001 / TestCP0.0001 columns in yellow ---
$wsRPart.Range("F1:G5").Interior.Color = 65535

# --- Update saved selections on both sheets ---
[void]$wsLinear.Range("C5").Select()
[void]$wsRPart.Range("E16").Select()
